# Update "Inscricoes" counts on the "Inscricoes" sheet of the Resumo de
# Inscricoes workbook. Each entry below lists the row number and the new
# values for columns E (Inscritos), F (Pagos) and H (Inscricoes
# homologadas) — column G (Isencoes deferidas) is untouched throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$updates = @(
    @{ Row = 8;  E = 7;   F = $null; H = $null },
    @{ Row = 15; E = 107; F = 45;    H = 45 },
    @{ Row = 17; E = 62;  F = $null; H = $null },
    @{ Row = 18; E = 58;  F = $null; H = $null },
    @{ Row = 26; E = 15;  F = $null; H = $null },
    @{ Row = 29; E = 8;   F = 3;     H = 3 },
    @{ Row = 37; E = 25;  F = 11;    H = 11 },
    @{ Row = 40; E = 8;   F = 3;     H = 3 },
    @{ Row = 43; E = 13;  F = $null; H = $null },
    @{ Row = 48; E = 14;  F = $null; H = $null },
    @{ Row = 49; E = 39;  F = $null; H = $null },
    @{ Row = 50; E = 12;  F = $null; H = $null },
    @{ Row = 53; E = 3;   F = $null; H = $null },
    @{ Row = 63; E = 14;  F = $null; H = $null },
    @{ Row = 71; E = 17;  F = $null; H = $null },
    @{ Row = 72; E = 21;  F = $null; H = $null },
    @{ Row = 88; E = 8;   F = 5;     H = 5 },
    @{ Row = 89; E = 17;  F = $null; H = $null }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("E$r").Value = $u.E
    if ($null -ne $u.F) { $ws.Range("F$r").Value = $u.F }
    if ($null -ne $u.H) { $ws.Range("H$r").Value = $u.H }
}
